# Updated cryptos list values (price/volume) per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.405.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -4.51%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.569.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.71%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -0.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'291.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.48%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3656"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.32%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'49.37"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.41%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -4.24%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.39%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07592"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -6.10%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D13").Value = "'21.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -5.10%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.897"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -6.04%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.00001141"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.99%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.568.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'89.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.06764"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.71%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'6.247"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -7.56%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.5297"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -7.91%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'16.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -5.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -3.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'22.416.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.54%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.395"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.06%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'3.008"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.05%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.73%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'144.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.81%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.965"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.49%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'125.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.60%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.738.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.045"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.277"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -9.61%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -7.64%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'10.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -9.48%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.02568"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.45%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.08439"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.61%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.2306"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -5.41%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.06540"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.66%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.543"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -6.54%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -8.47%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.256"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.6397"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -7.08%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'14.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.82%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.00%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.6035"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.15%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.784"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.18%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.135"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -5.12%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'122.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.76%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +2.54%  "
$ws.Range("E51").Style = "Normal"
